$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value = 47.23036199999999
$ws.Range("H2").Value = 141.691086
$ws.Range("I2").Value = 0.3244251370417807
$ws.Range("J2").Value = 0.3244251370417807
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.03057133333333333
$ws.Range("N2").Value = 0.09171399999999999
$ws.Range("O2").Value = 0.1084248755136686
$ws.Range("P2").Value = 0.1084248755136687
$ws.Range("Q2").Value = 1.443895140156
$ws.Range("R2").Value = 12.995056261404
$ws.Range("S2").Value = 0.03517575509725996
$ws.Range("T2").Value = 0.03517575509725997

# Row 3
$ws.Range("G3").Value = 47.23036199999999
$ws.Range("H3").Value = 141.691086
$ws.Range("I3").Value = 0.3244251370417807
$ws.Range("J3").Value = 0.3244251370417807
$ws.Range("O3").Value = 0.8915751244863314
$ws.Range("P3").Value = 0.8915751244863314
$ws.Range("Q3").Value = 11.873114755548
$ws.Range("R3").Value = 106.858032799932
$ws.Range("S3").Value = 0.2892493819445208
$ws.Range("T3").Value = 0.2892493819445208

# Row 4
$ws.Range("I4").Value = 0.4188548944674916
$ws.Range("J4").Value = 0.4188548944674916
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.03057133333333333
$ws.Range("N4").Value = 0.09171399999999999
$ws.Range("O4").Value = 0.1084248755136686
$ws.Range("P4").Value = 0.1084248755136687
$ws.Range("Q4").Value = 1.864166729085111
$ws.Range("R4").Value = 16.777500561766
$ws.Range("S4").Value = 0.04541428979092859
$ws.Range("T4").Value = 0.0454142897909286

# Row 5
$ws.Range("I5").Value = 0.4188548944674916
$ws.Range("J5").Value = 0.4188548944674916
$ws.Range("O5").Value = 0.8915751244863314
$ws.Range("P5").Value = 0.8915751244863314
$ws.Range("S5").Value = 0.373440604676563
$ws.Range("T5").Value = 0.373440604676563

# Row 6
$ws.Range("I6").Value = 0.2567199684907278
$ws.Range("J6").Value = 0.2567199684907277
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.03057133333333333
$ws.Range("N6").Value = 0.09171399999999999
$ws.Range("O6").Value = 0.1084248755136686
$ws.Range("P6").Value = 0.1084248755136687
$ws.Range("Q6").Value = 1.142564716978222
$ws.Range("R6").Value = 10.283082452804
$ws.Range("S6").Value = 0.02783483062548009
$ws.Range("T6").Value = 0.02783483062548009

# Row 7
$ws.Range("I7").Value = 0.2567199684907278
$ws.Range("J7").Value = 0.2567199684907277
$ws.Range("O7").Value = 0.8915751244863314
$ws.Range("P7").Value = 0.8915751244863314
$ws.Range("S7").Value = 0.2288851378652477
$ws.Range("T7").Value = 0.2288851378652476
